# Auto-generated edit script: rebuild sheet1 data per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: add new header count column T1
$ws.Range("T1").Value = 18
$ws.Range("S1").Copy()
$ws.Range("T1").PasteSpecial(-4122)

# Row 2: HKL header labels shifted right by one (new "1Pair-B" col) + new MaxUnique col T2
$ws.Range("B2").Value = "HKL"
$ws.Range("C2").Value = "[1, 1, 0]"
$ws.Range("D2").Value = "[2, 0, 0]"
$ws.Range("E2").Value = "[2, 1, 1]"
$ws.Range("F2").Value = "[2, 2, 0]"
$ws.Range("G2").Value = "[3, 1, 0]"
$ws.Range("H2").Value = "[2, 2, 2]"
$ws.Range("I2").Value = "[3, 2, 1]"
$ws.Range("J2").Value = "[4, 0, 0]"
$ws.Range("K2").Value = "1Pair-A"
$ws.Range("L2").Value = "1Pair-B"
$ws.Range("M2").Value = "2Pairs-A"
$ws.Range("N2").Value = "2Pairs-B"
$ws.Range("O2").Value = "3Pairs-A"
$ws.Range("P2").Value = "3Pairs-B"
$ws.Range("Q2").Value = "3Pairs-C"
$ws.Range("R2").Value = "4Pairs"
$ws.Range("S2").Value = "5A4F"
$ws.Range("T2").Value = "MaxUnique"

# Rows 3-7: updated numeric data (existing rows) + new column T
# Row 3
$ws.Range("B3").Value = "Equal Angle"
$ws.Range("C3").Value = 1.110965417867435
$ws.Range("D3").Value = 0.8732204610951009
$ws.Range("E3").Value = 0.9870317002881844
$ws.Range("F3").Value = 1.110965417867435
$ws.Range("G3").Value = 0.9369596541786743
$ws.Range("H3").Value = 0.9840417867435158
$ws.Range("I3").Value = 1.017406340057637
$ws.Range("J3").Value = 0.8732204610951009
$ws.Range("K3").Value = 1.110965417867435
$ws.Range("L3").Value = 0.9870317002881844
$ws.Range("M3").Value = 0.9301260806916427
$ws.Range("N3").Value = 0.9301260806916427
$ws.Range("O3").Value = 0.9324039385206532
$ws.Range("P3").Value = 0.9904058597502402
$ws.Range("Q3").Value = 0.9904058597502402
$ws.Range("R3").Value = 1.020545749279539
$ws.Range("S3").Value = 1.020545749279539
$ws.Range("T3").Value = 0.9849375600384246

# Row 4
$ws.Range("B4").Value = "CLR"
$ws.Range("C4").Value = 0.9870975407402144
$ws.Range("D4").Value = 0.9950939054233108
$ws.Range("E4").Value = 0.9991745937072521
$ws.Range("F4").Value = 0.9870975407402144
$ws.Range("G4").Value = 0.9884895757715156
$ws.Range("H4").Value = 1.014622979784832
$ws.Range("I4").Value = 0.9968633548783302
$ws.Range("J4").Value = 0.9950939054233108
$ws.Range("K4").Value = 0.9870975407402144
$ws.Range("L4").Value = 0.9991745937072521
$ws.Range("M4").Value = 0.9971342495652815
$ws.Range("N4").Value = 0.9971342495652815
$ws.Range("O4").Value = 0.9942526916340263
$ws.Range("P4").Value = 0.9937886799569258
$ws.Range("Q4").Value = 0.9937886799569258
$ws.Range("R4").Value = 0.9921158951527479
$ws.Range("S4").Value = 0.9921158951527479
$ws.Range("T4").Value = 0.9968903250509092

# Row 5
$ws.Range("B5").Value = "BT8Hex"
$ws.Range("C5").Value = 0.9833508948227769
$ws.Range("D5").Value = 0.9950887042466666
$ws.Range("E5").Value = 1.002219417062731
$ws.Range("F5").Value = 0.9833508948227769
$ws.Range("G5").Value = 0.984872430553372
$ws.Range("H5").Value = 1.026991197513424
$ws.Range("I5").Value = 0.9979006955028106
$ws.Range("J5").Value = 0.9950887042466666
$ws.Range("K5").Value = 0.9833508948227769
$ws.Range("L5").Value = 1.002219417062731
$ws.Range("M5").Value = 0.9986540606546985
$ws.Range("N5").Value = 0.9986540606546985
$ws.Range("O5").Value = 0.9940601839542564
$ws.Range("P5").Value = 0.9935530053773913
$ws.Range("Q5").Value = 0.9935530053773913
$ws.Range("R5").Value = 0.9910024777387377
$ws.Range("S5").Value = 0.9910024777387377
$ws.Range("T5").Value = 0.9984038899502967

# Row 6
$ws.Range("B6").Value = "Spiral"
$ws.Range("C6").Value = 0.9962784887745791
$ws.Range("D6").Value = 0.994619251703808
$ws.Range("E6").Value = 0.9936797700849543
$ws.Range("F6").Value = 0.9962784887745791
$ws.Range("G6").Value = 0.9933988029901977
$ws.Range("H6").Value = 0.9941638581492266
$ws.Range("I6").Value = 0.9954614488714957
$ws.Range("J6").Value = 0.994619251703808
$ws.Range("K6").Value = 0.9962784887745791
$ws.Range("L6").Value = 0.9936797700849543
$ws.Range("M6").Value = 0.9941495108943812
$ws.Range("N6").Value = 0.9941495108943812
$ws.Range("O6").Value = 0.9938992749263201
$ws.Range("P6").Value = 0.9948591701877806
$ws.Range("Q6").Value = 0.9948591701877806
$ws.Range("R6").Value = 0.9952139998344802
$ws.Range("S6").Value = 0.9952139998344802
$ws.Range("T6").Value = 0.9946002700957104

# Row 7
$ws.Range("B7").Value = "OffsetF"
$ws.Range("C7").Value = 0.8944693780818429
$ws.Range("D7").Value = 0.1435974566754628
$ws.Range("E7").Value = 1.398067464417403
$ws.Range("F7").Value = 0.8944693780818429
$ws.Range("G7").Value = 0.3856936744933193
$ws.Range("H7").Value = 2.243664783995547
$ws.Range("I7").Value = 1.253513552519149
$ws.Range("J7").Value = 0.1435974566754628
$ws.Range("K7").Value = 0.8944693780818429
$ws.Range("L7").Value = 1.398067464417403
$ws.Range("M7").Value = 0.7708324605464331
$ws.Range("N7").Value = 0.7708324605464331
$ws.Range("O7").Value = 0.6424528651953951
$ws.Range("P7").Value = 0.8120447663915696
$ws.Range("Q7").Value = 0.8120447663915696
$ws.Range("R7").Value = 0.832650919314138
$ws.Range("S7").Value = 0.832650919314138
$ws.Range("T7").Value = 1.053167718363788

# Rows 8-11: new rows inserted (OffsetA, RD Single, TD Single, HexGrid-90degTilt5degRes)
# Row 8
$ws.Range("A8").Value = 6
$ws.Range("A2").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("B8").Value = "OffsetA"
$ws.Range("C8").Value = 1.22692954780446
$ws.Range("D8").Value = 0.5449080945638909
$ws.Range("E8").Value = 1.041304524754317
$ws.Range("F8").Value = 1.22692954780446
$ws.Range("G8").Value = 0.8145427962352142
$ws.Range("H8").Value = 0.9800356093511665
$ws.Range("I8").Value = 1.094669397544617
$ws.Range("J8").Value = 0.5449080945638909
$ws.Range("K8").Value = 1.22692954780446
$ws.Range("L8").Value = 1.041304524754317
$ws.Range("M8").Value = 0.793106309659104
$ws.Range("N8").Value = 0.793106309659104
$ws.Range("O8").Value = 0.800251805184474
$ws.Range("P8").Value = 0.9377140557075562
$ws.Range("Q8").Value = 0.9377140557075562
$ws.Range("R8").Value = 1.010017928731782
$ws.Range("S8").Value = 1.010017928731782
$ws.Range("T8").Value = 0.9503983283756111

# Row 9
$ws.Range("A9").Value = 7
$ws.Range("A2").Copy()
$ws.Range("A9").PasteSpecial(-4122)
$ws.Range("B9").Value = "RD Single"
$ws.Range("C9").Value = 0.21
$ws.Range("D9").Value = 3.94
$ws.Range("E9").Value = 0.46
$ws.Range("F9").Value = 0.21
$ws.Range("G9").Value = 2.23
$ws.Range("H9").Value = 0.07
$ws.Range("I9").Value = 0.4
$ws.Range("J9").Value = 3.94
$ws.Range("K9").Value = 0.21
$ws.Range("L9").Value = 0.46
$ws.Range("M9").Value = 2.2
$ws.Range("N9").Value = 2.2
$ws.Range("O9").Value = 2.21
$ws.Range("P9").Value = 1.536666666666667
$ws.Range("Q9").Value = 1.536666666666667
$ws.Range("R9").Value = 1.205
$ws.Range("S9").Value = 1.205
$ws.Range("T9").Value = 1.218333333333333

# Row 10
$ws.Range("A10").Value = 8
$ws.Range("A2").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("B10").Value = "TD Single"
$ws.Range("C10").Value = 1.97
$ws.Range("D10").Value = 0.21
$ws.Range("E10").Value = 0.86
$ws.Range("F10").Value = 1.97
$ws.Range("G10").Value = 0.64
$ws.Range("H10").Value = 0.67
$ws.Range("I10").Value = 1.11
$ws.Range("J10").Value = 0.21
$ws.Range("K10").Value = 1.97
$ws.Range("L10").Value = 0.86
$ws.Range("M10").Value = 0.535
$ws.Range("N10").Value = 0.535
$ws.Range("O10").Value = 0.57
$ws.Range("P10").Value = 1.013333333333333
$ws.Range("Q10").Value = 1.013333333333333
$ws.Range("R10").Value = 1.2525
$ws.Range("S10").Value = 1.2525
$ws.Range("T10").Value = 0.9100000000000001

# Row 11
$ws.Range("A11").Value = 9
$ws.Range("A2").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("B11").Value = "HexGrid-90degTilt5degRes"
$ws.Range("C11").Value = 0.9915490378394113
$ws.Range("D11").Value = 0.9974757325517742
$ws.Range("E11").Value = 0.994935224634625
$ws.Range("F11").Value = 0.9915490378394113
$ws.Range("G11").Value = 0.9937214692675462
$ws.Range("H11").Value = 0.9990330497856222
$ws.Range("I11").Value = 0.9952945993521713
$ws.Range("J11").Value = 0.9974757325517742
$ws.Range("K11").Value = 0.9915490378394113
$ws.Range("L11").Value = 0.994935224634625
$ws.Range("M11").Value = 0.9962054785931996
$ws.Range("N11").Value = 0.9962054785931996
$ws.Range("O11").Value = 0.9953774754846485
$ws.Range("P11").Value = 0.9946533316752703
$ws.Range("Q11").Value = 0.9946533316752703
$ws.Range("R11").Value = 0.9938772582163055
$ws.Range("S11").Value = 0.9938772582163055
$ws.Range("T11").Value = 0.9953348522385251

